$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("userlive")

# Order matters for shared-string table indices: write in the same
# sequence the target workbook's <sst> grew (A2@Sheet1, B2@userlive,
# F2@Sheet1, B2@Sheet1).
$ws1.Range("A2").Value = "minhson0907"
$ws2.Range("B2").Value = "'"
$ws1.Range("F2").Value = "aGlnaGxhbmQxMFg="
$ws1.Range("B2").Value = "'002704070016025"

# Sheet2 selection/tab moves from E15 to C12, and tabSelected is dropped.
$ws2.Select()
$ws2.Range("C12").Select()

# Sheet1 becomes the selected/active tab with selection on G18.
$ws1.Select()
$ws1.Range("G18").Select()
